$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# manualStatus (column I) values for rows 8-10 change from the number 4
# to the text "[4]"
$ws.Range("I8").Value = "[4]"
$ws.Range("I9").Value = "[4]"
$ws.Range("I10").Value = "[4]"

# Rows 9 and 10 got a little shorter once the new text went in
$ws.Rows.Item(9).RowHeight = 13.8
$ws.Rows.Item(10).RowHeight = 13.8

# Column F (fastqFileName) is widened so the long file names are visible
$ws.Columns.Item(6).ColumnWidth = 48

# Selection moved to I10
[void]$ws.Range("I10").Select()
